# Upravene l0m jako lm pri PP
# Replace the values in column C ("Optimal_length_upravene") with the
# adjusted optimal length values (l0m as lm at PP), keeping the cells as
# text (shared-string) entries rather than numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$refStyle = $ws.Range("A1").Style

$values = @(
  "0.0693861569985691",
  "0.116982312421667",
  "0.142349225195613",
  "0.169624773492201",
  "0.142003456015892",
  "0.190963063904747",
  "0.122063988214140",
  "0.147660576545422",
  "0.0874282936631438",
  "0.0755448185556480",
  "0.0638454741557207",
  "0.0813852241566054",
  "0.0808168514025335"
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $cell = $ws.Cells.Item($row, 3)
  $cell.NumberFormat = "@"
  $cell.Value = $values[$i]
  $cell.Style = $refStyle
}

# Move the active selection to C14, matching the author's final cursor position
$ws.Range("C14").Select()
